# Update "想去人数" (F) counts and the sold-out marker in G6 on both the
# "展览" and "全部类型" sheets (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1783
    3  = 251
    4  = 230
    5  = 7614
    6  = 571
    7  = 527
    8  = 75
    9  = 20
    10 = 9085
    11 = 2382
    12 = 295
    13 = 9841
    14 = 10486
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }

    # Row 6 sold out -> replace numeric price with the "sold out" label.
    $ws.Cells.Item(6, 7).Value = "已售罄"
}
